$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D/E (and B/C) cells in this sheet are free-form text (coin price /
# volume strings), never real numbers. Force text format on any target
# cell whose new content would otherwise be auto-coerced to a number by
# Excel (e.g. "1.00", "428.58") so it round-trips as the literal string.

$ws.Range('D2').Value = '66.115.38'
$ws.Range('E2').Value = '  +3.46%  '

$ws.Range('D3').Value = '3.833.39'
$ws.Range('E3').Value = '  +7.50%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '428.58'
$ws.Range('E5').Value = '  +8.80%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.75'
$ws.Range('E6').Value = '  +4.27%  '

$ws.Range('D7').Value = '3.828.53'
$ws.Range('E7').Value = '  +7.76%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.614'
$ws.Range('E8').Value = '  +3.68%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.10%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.740'
$ws.Range('E10').Value = '  +7.63%  '

$ws.Range('E11').Value = '  +3.60%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000338'
$ws.Range('E12').Value = '  -2.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.12'
$ws.Range('E13').Value = '  +7.35%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.71'
$ws.Range('E14').Value = '  +14.87%  '

$ws.Range('D15').Value = '4.434.55'
$ws.Range('E15').Value = '  +7.37%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.43'
$ws.Range('E16').Value = '  +18.20%  '

$ws.Range('E17').Value = '  +0.69%  '

$ws.Range('D18').Value = '3.838.20'
$ws.Range('E18').Value = '  +8.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.16'
$ws.Range('E19').Value = '  +6.59%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.13'
$ws.Range('E20').Value = '  +9.28%  '

$ws.Range('D21').Value = '66.414.19'
$ws.Range('E21').Value = '  +3.88%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '417.48'
$ws.Range('E22').Value = '  +4.55%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.57'
$ws.Range('E23').Value = '  +11.98%  '

$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.19'
$ws.Range('E24').Value = '  +10.96%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.18'
$ws.Range('E25').Value = '  +5.36%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.43'
$ws.Range('E26').Value = '  +9.53%  '

$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.31'
$ws.Range('E27').Value = '  +10.07%  '

$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('E28').Value = '  +13.44%  '

$ws.Range('E29').Value = '  -0.29%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.33'
$ws.Range('E30').Value = '  +35.80%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '14.48'
$ws.Range('E31').Value = '  +21.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '715.59'
$ws.Range('E32').Value = '  +6.62%  '

$ws.Range('E33').Value = '  +13.51%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.78'
$ws.Range('E34').Value = '  +4.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.95'
$ws.Range('E35').Value = '  +42.33%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.14%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.154'
$ws.Range('E37').Value = '  +1.99%  '

$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '39.27'
$ws.Range('E38').Value = '  +5.07%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '56.23'
$ws.Range('E39').Value = '  +3.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0473'
$ws.Range('E40').Value = '  +5.20%  '

$ws.Range('D41').Value = '0.0₃0712'
$ws.Range('E41').Value = '  +15.46%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.93'
$ws.Range('E42').Value = '  +4.09%  '

$ws.Range('E43').Value = '  +0.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.139'
$ws.Range('E44').Value = '  +4.68%  '

$ws.Range('B45').Value = 'LidoDAOToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.42'
$ws.Range('E45').Value = '  +9.77%  '

$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.23'
$ws.Range('E46').Value = '  +6.88%  '

$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.47'
$ws.Range('E47').Value = '  +45.20%  '

$ws.Range('B48').Value = 'TheGraph'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.322'
$ws.Range('E48').Value = '  +16.20%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.08'
$ws.Range('E49').Value = '  +5.15%  '

$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.64'
$ws.Range('E50').Value = '  +5.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '141.67'
